# Modify the phenotyep-microbiome - vaginal dryness & Endomentris
# Column E = "Endometritis" probability, Column K = "Vaginal Dryness" probability
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 49.1
$ws.Range("K2").Value2 = 53.2
$ws.Range("K3").Value2 = 5
$ws.Range("E4").Value2 = 44.5
$ws.Range("K4").Value2 = 81.5
$ws.Range("E5").Value2 = 54
$ws.Range("K5").Value2 = 30
$ws.Range("K6").Value2 = 73.5
$ws.Range("E7").Value2 = 8.300000000000001
$ws.Range("K7").Value2 = 29.7
$ws.Range("E8").Value2 = 8.6
$ws.Range("K8").Value2 = 74.59999999999999
$ws.Range("E9").Value2 = 21.5
$ws.Range("K9").Value2 = 23.7
$ws.Range("E10").Value2 = 72
$ws.Range("K10").Value2 = 54.1
$ws.Range("E11").Value2 = 9.1
$ws.Range("K11").Value2 = 89.7
$ws.Range("E12").Value2 = 58.4
$ws.Range("K12").Value2 = 88.3
$ws.Range("E13").Value2 = 33
$ws.Range("K13").Value2 = 81.7
$ws.Range("E14").Value2 = 77.7
$ws.Range("K14").Value2 = 94.40000000000001
$ws.Range("E15").Value2 = 71.40000000000001
$ws.Range("K15").Value2 = 91.2
$ws.Range("E16").Value2 = 72.90000000000001
$ws.Range("K16").Value2 = 16.7
$ws.Range("E17").Value2 = 55.6
$ws.Range("K17").Value2 = 42.4
$ws.Range("E18").Value2 = 74.3
$ws.Range("K18").Value2 = 83.2
$ws.Range("E19").Value2 = 75.40000000000001
$ws.Range("K19").Value2 = 66.2
$ws.Range("E20").Value2 = 67.5
$ws.Range("K20").Value2 = 13.6
$ws.Range("K21").Value2 = 68
$ws.Range("E22").Value2 = 80.59999999999999
$ws.Range("K22").Value2 = 51
$ws.Range("E23").Value2 = 92.8
$ws.Range("K23").Value2 = 92.2
$ws.Range("E24").Value2 = 60.3
$ws.Range("K24").Value2 = 68
$ws.Range("E25").Value2 = 15
$ws.Range("K25").Value2 = 5
$ws.Range("E26").Value2 = 57.7
$ws.Range("K26").Value2 = 28
$ws.Range("E27").Value2 = 7.6
$ws.Range("K27").Value2 = 73.40000000000001
$ws.Range("E28").Value2 = 43.8
$ws.Range("K28").Value2 = 66.3
$ws.Range("E29").Value2 = 66.59999999999999
$ws.Range("K29").Value2 = 85.59999999999999
$ws.Range("E30").Value2 = 75.59999999999999
$ws.Range("K30").Value2 = 91.09999999999999
$ws.Range("K31").Value2 = 60.8
$ws.Range("K32").Value2 = 58.2
$ws.Range("E33").Value2 = 40
$ws.Range("K33").Value2 = 5
$ws.Range("E34").Value2 = 37.9
$ws.Range("K34").Value2 = 15
$ws.Range("E35").Value2 = 93
$ws.Range("K35").Value2 = 91.7
$ws.Range("K36").Value2 = 15.1
$ws.Range("E37").Value2 = 54.1
$ws.Range("K37").Value2 = 50.3
$ws.Range("E38").Value2 = 44.2
$ws.Range("K38").Value2 = 54.4
$ws.Range("E39").Value2 = 63.2
$ws.Range("K39").Value2 = 46
$ws.Range("K40").Value2 = 48.4
$ws.Range("E41").Value2 = 35.8
$ws.Range("K41").Value2 = 57.5
$ws.Range("E42").Value2 = 69
$ws.Range("K42").Value2 = 79.7
$ws.Range("E43").Value2 = 7.2
$ws.Range("K43").Value2 = 14.8
$ws.Range("E44").Value2 = 72.09999999999999
$ws.Range("K44").Value2 = 88
$ws.Range("E45").Value2 = 53.9
$ws.Range("K45").Value2 = 88.2
$ws.Range("E46").Value2 = 68.2
$ws.Range("K46").Value2 = 27.8
$ws.Range("E47").Value2 = 16.3
$ws.Range("K47").Value2 = 20.9
$ws.Range("E48").Value2 = 50.9
$ws.Range("K48").Value2 = 5
$ws.Range("E49").Value2 = 41
$ws.Range("K49").Value2 = 14.4
$ws.Range("E50").Value2 = 38.5
$ws.Range("K50").Value2 = 59.3
$ws.Range("E51").Value2 = 6.6
$ws.Range("K51").Value2 = 39.5
$ws.Range("E52").Value2 = 37.5
$ws.Range("K52").Value2 = 14.9
$ws.Range("E53").Value2 = 61.3
$ws.Range("K53").Value2 = 76.2
$ws.Range("E54").Value2 = 45.6
$ws.Range("K54").Value2 = 5
$ws.Range("E55").Value2 = 63.7
$ws.Range("K55").Value2 = 47.4
$ws.Range("E56").Value2 = 71.90000000000001
$ws.Range("K56").Value2 = 17
$ws.Range("K57").Value2 = 48.7
$ws.Range("K58").Value2 = 90.3
$ws.Range("K59").Value2 = 79.09999999999999
$ws.Range("E60").Value2 = 50.1
$ws.Range("K60").Value2 = 54
$ws.Range("E61").Value2 = 42.5
$ws.Range("K61").Value2 = 15.5
$ws.Range("E62").Value2 = 24.5
$ws.Range("K62").Value2 = 41.5
$ws.Range("K63").Value2 = 5
$ws.Range("E64").Value2 = 39.8
$ws.Range("K64").Value2 = 20.9
$ws.Range("E65").Value2 = 30.5
$ws.Range("K65").Value2 = 13.8
$ws.Range("E66").Value2 = 38.3
$ws.Range("K66").Value2 = 59
$ws.Range("E67").Value2 = 5
$ws.Range("K67").Value2 = 14.2
$ws.Range("E68").Value2 = 23.6
$ws.Range("K68").Value2 = 61.9
$ws.Range("E69").Value2 = 29.9
$ws.Range("K69").Value2 = 18.7
$ws.Range("E70").Value2 = 73.7
$ws.Range("K70").Value2 = 28.8
$ws.Range("E71").Value2 = 45.7
$ws.Range("K71").Value2 = 91.90000000000001
$ws.Range("E72").Value2 = 51
$ws.Range("K72").Value2 = 82
$ws.Range("K73").Value2 = 15.2
$ws.Range("E74").Value2 = 18.3
$ws.Range("K74").Value2 = 17.9
$ws.Range("E75").Value2 = 33.1
$ws.Range("K75").Value2 = 77.09999999999999
$ws.Range("E77").Value2 = 82.5
$ws.Range("K77").Value2 = 17.2
$ws.Range("E78").Value2 = 30.8
$ws.Range("K78").Value2 = 78.3
$ws.Range("E79").Value2 = 17.3
$ws.Range("K79").Value2 = 13.6
$ws.Range("E80").Value2 = 64.59999999999999
$ws.Range("K80").Value2 = 46.3
$ws.Range("E81").Value2 = 61.9
$ws.Range("K81").Value2 = 71.8
$ws.Range("E82").Value2 = 58.1
$ws.Range("K82").Value2 = 55.7
$ws.Range("E83").Value2 = 92.2
$ws.Range("E84").Value2 = 38.4
$ws.Range("K84").Value2 = 72.5
$ws.Range("E85").Value2 = 83
$ws.Range("K85").Value2 = 60.4
$ws.Range("E86").Value2 = 24.4
$ws.Range("K86").Value2 = 41.4
$ws.Range("E87").Value2 = 30
$ws.Range("K87").Value2 = 74.90000000000001
$ws.Range("E88").Value2 = 71
$ws.Range("K88").Value2 = 69.7
$ws.Range("E89").Value2 = 72.90000000000001
$ws.Range("K89").Value2 = 58.1
$ws.Range("E90").Value2 = 25.9
$ws.Range("K90").Value2 = 40.5
$ws.Range("E91").Value2 = 43.5
$ws.Range("K91").Value2 = 5
$ws.Range("E92").Value2 = 49.1
$ws.Range("K92").Value2 = 41.9
